# Update the "thresholds_summary" results table and the chosen-threshold
# callouts to reflect the refreshed pipeline metrics.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Columns: 1=threshold 2=row_count 3=unique_count 4=repeated_count 5=no_repeats_bool

# threshold = 0.6  (row 14)
$t.Cell(14, 2).Range.Text = "6"
$t.Cell(14, 3).Range.Text = "6"

# threshold = 0.65 (row 15)
$t.Cell(15, 2).Range.Text = "7"
$t.Cell(15, 3).Range.Text = "7"

# threshold = 0.7  (row 16)
$t.Cell(16, 2).Range.Text = "7"
$t.Cell(16, 3).Range.Text = "7"

# threshold = 0.75 (row 17)
$t.Cell(17, 2).Range.Text = "8"
$t.Cell(17, 3).Range.Text = "8"

# threshold = 0.8  (row 18)
$t.Cell(18, 2).Range.Text = "10"
$t.Cell(18, 3).Range.Text = "10"

# threshold = 0.95 (row 21)
$t.Cell(21, 2).Range.Text = "11"
$t.Cell(21, 4).Range.Text = "2"
$t.Cell(21, 5).Range.Text = "False"

# threshold = 1.0  (row 22)
$t.Cell(22, 2).Range.Text = "11"
$t.Cell(22, 4).Range.Text = "2"
$t.Cell(22, 5).Range.Text = "False"

# Update chosen-threshold summary paragraphs.
$d.Content.Find.Execute("Full Analysis threshold = 1.00.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Full Analysis threshold = 0.90.", 2)
$d.Content.Find.Execute("Core-Level Analysis threshold = 0.80.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Core-Level Analysis threshold = 0.75.", 2)
